# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (cloning the "2021-Q4" sheet's layout/style)
# right after "2021-Q4" and before "总计", populates it with the single
# fund-holding row, and updates the "总计" (summary) sheet with a new
# leading row for 2022-Q1, shifting the existing rows down and
# renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Clone header/formatting from the "2021-Q4" sheet (same column layout:
# 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名)
$afterSheet.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$afterSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header text
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Data row (fund holding)
$newSheet.Cells.Item(2, 1).Value = 0

$c = $newSheet.Cells.Item(2, 2)
$c.Value = "'001118"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2, 3)
$c.Value = "华宝事件驱动混合"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2, 4)
$c.Value = "'6.29"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2, 5)
$c.Value = "'92.73"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2, 6)
$c.Value = "'3.11"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2, 7)
$c.Value = "'0.1956"
$c.Style = "Normal"

$newSheet.Cells.Item(2, 8).Value = 9

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q1
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows(2).Insert()
$summary.Range("B2:D2").Style = "Normal"

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q1"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.2

# Restore the index-column style on the new row (matches the other rows)
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# Renumber the index column (A) for the rows pushed down
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
